# Update to final publication zib MedicalDevice
#
# The sheet gains two new columns - "alias_zib" (E) and "stereotype_zib" (H) -
# inserted into the existing zib/xtehr mapping table, a stray duplicate row is
# removed, and a handful of "definitioncode_zib" cells get a coding-system
# prefix (e.g. "SNOMED CT: ...").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the two new columns at their final positions.
#    Inserting "E" first shifts the old E..P block to F..Q, then inserting
#    "H" shifts the old G..P (now at H..Q) block to I..R - leaving E and H
#    empty and ready for the new alias_zib / stereotype_zib values, and every
#    other former column landing exactly where the new header row expects it.
$ws.Columns("E").Insert()
$ws.Columns("H").Insert()

# 2) Drop the stray duplicate row (old row 22: EHDSDeviceUse.reason /
#    MedicalDevice.Indication::Problem) - everything below shifts up by one.
$ws.Rows(22).Delete()

# 3) Fill in the new header labels.
$ws.Range("E1").Value = "alias_zib"
$ws.Range("H1").Value = "stereotype_zib"

# 4) Fill in the new alias_zib / stereotype_zib values for the rows that have
#    zib content.
$ws.Range("E2").Value = "NL: Product"
$ws.Range("H2").Value = "container"

$ws.Range("E4").Value = "NL: ProductID"
$ws.Range("H4").Value = "data"

$ws.Range("E12").Value = "NL: ProductType"
$ws.Range("H12").Value = "data"

$ws.Range("E14").Value = "NL: MedischHulpmiddel"
$ws.Range("H14").Value = "rootconcept"

$ws.Range("E15").Value = "NL: AnatomischeLocatie"
$ws.Range("H15").Value = "data,reference"

$ws.Range("E17").Value = "NL: EindDatum"
$ws.Range("H17").Value = "data"

$ws.Range("E19").Value = "NL: BeginDatum"
$ws.Range("H19").Value = "data"

$ws.Range("E20").Value = "NL: Toelichting"
$ws.Range("H20").Value = "data"

$ws.Range("E21").Value = "NL: Indicatie::Diagnose"
$ws.Range("H21").Value = "context,reference"

$ws.Range("E28").Value = "NL: Zorgverlener"
$ws.Range("H28").Value = "context,reference"

$ws.Range("E29").Value = "NL: Locatie::Zorgaanbieder"
$ws.Range("H29").Value = "context,reference"

$ws.Range("E30").Value = "NL: ProductOmschrijving"
$ws.Range("H30").Value = "data"

# 5) Prefix the definitioncode_zib (now column K) values with their coding
#    system name.
$ws.Range("K2").Value = "SNOMED CT: 405815000 Procedure device"
$ws.Range("K14").Value = "SNOMED CT: 49062001 Device"
$ws.Range("K15").Value = "SNOMED CT: 363698007 Finding site"
$ws.Range("K20").Value = "LOINC: 48767-8 Annotation comment [Interpretation] Narrative"
